$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps its original text storage; values such as
# "1.00" or "239.54" would otherwise be auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '91.515.06'
$ws.Range("E2").Value = '  +1.07%  '
$ws.Range("D3").Value = '3.140.11'
$ws.Range("E3").Value = '  +1.31%  '
$ws.Range("E4").Value = '  +0.17%  '
$ws.Range("D5").Value = '239.54'
$ws.Range("E5").Value = '  -0.98%  '
$ws.Range("D6").Value = '617.88'
$ws.Range("E6").Value = '  -1.18%  '
$ws.Range("D7").Value = '1.11'
$ws.Range("E7").Value = '  -2.13%  '
$ws.Range("D8").Value = '0.386'
$ws.Range("E8").Value = '  +3.17%  '
$ws.Range("E9").Value = '  +0.01%  '
$ws.Range("D10").Value = '3.147.79'
$ws.Range("E10").Value = '  +14.40%  '
$ws.Range("D11").Value = '0.742'
$ws.Range("E11").Value = '  +0.48%  '
$ws.Range("E12").Value = '  +0.34%  '
$ws.Range("E13").Value = '  +1.12%  '
$ws.Range("D14").Value = '34.90'
$ws.Range("E14").Value = '  -0.89%  '
$ws.Range("D15").Value = '5.58'
$ws.Range("E15").Value = '  +1.82%  '
$ws.Range("D16").Value = '91.280.53'
$ws.Range("E16").Value = '  +0.93%  '
$ws.Range("E17").Value = '  +1.55%  '
$ws.Range("D18").Value = '3.153.77'
$ws.Range("E18").Value = '  -0.27%  '
$ws.Range("E19").Value = '  -3.13%  '
$ws.Range("D20").Value = '14.91'
$ws.Range("E20").Value = '  +4.82%  '
$ws.Range("D21").Value = '5.85'
$ws.Range("E21").Value = '  +1.74%  '
$ws.Range("D22").Value = '455.57'
$ws.Range("E22").Value = '  +2.22%  '
$ws.Range("D23").Value = '0.0000202'
$ws.Range("E23").Value = '  -3.79%  '
$ws.Range("D24").Value = '9.16'
$ws.Range("E24").Value = '  +1.00%  '
$ws.Range("D25").Value = '5.86'
$ws.Range("E25").Value = '  -0.15%  '
$ws.Range("E26").Value = '  +62.04%  '
$ws.Range("D27").Value = '88.50'
$ws.Range("E27").Value = '  -4.88%  '
$ws.Range("D28").Value = '11.74'
$ws.Range("E28").Value = '  -2.32%  '
$ws.Range("B29").Value = 'Hedera'
$ws.Range("C29").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D29").Value = '0.152'
$ws.Range("E29").Value = '  +42.60%  '
$ws.Range("B30").Value = 'WrappedeETH'
$ws.Range("C30").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D30").Value = '3.317.67'
$ws.Range("E30").Value = '  +1.76%  '
$ws.Range("D32").Value = '0.228'
$ws.Range("E32").Value = '  +4.57%  '
$ws.Range("D33").Value = '0.166'
$ws.Range("E33").Value = '  -4.97%  '
$ws.Range("D34").Value = '9.36'
$ws.Range("E34").Value = '  +1.27%  '
$ws.Range("D35").Value = '0.175'
$ws.Range("E35").Value = '  +11.38%  '
$ws.Range("D36").Value = '26.25'
$ws.Range("E36").Value = '  -1.10%  '
$ws.Range("D37").Value = '7.44'
$ws.Range("E37").Value = '  -0.99%  '
$ws.Range("D38").Value = '2.01'
$ws.Range("E38").Value = '  +4.79%  '
$ws.Range("B39").Value = 'MantraDAO'
$ws.Range("C39").Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range("D39").Value = '3.92'
$ws.Range("E39").Value = '  -10.58%  '
$ws.Range("B40").Value = 'Bittensor'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D40").Value = '491.26'
$ws.Range("E40").Value = '  -0.31%  '
$ws.Range("E41").Value = '  +2.20%  '
$ws.Range("D42").Value = '0.440'
$ws.Range("E42").Value = '  +5.65%  '
$ws.Range("E43").Value = '  -6.02%  '
$ws.Range("D44").Value = '22.14'
$ws.Range("E44").Value = '  +0.10%  '
$ws.Range("E45").Value = '  +0.00%  '
$ws.Range("B46").Value = 'ARBITRUM'
$ws.Range("C46").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D46").Value = '0.705'
$ws.Range("E46").Value = '  +3.06%  '
$ws.Range("B47").Value = 'Stacks'
$ws.Range("C47").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D47").Value = '1.92'
$ws.Range("E47").Value = '  +1.29%  '
$ws.Range("D48").Value = '155.71'
$ws.Range("E48").Value = '  -2.23%  '
$ws.Range("D49").Value = '1.35'
$ws.Range("E49").Value = '  +1.41%  '
$ws.Range("D50").Value = '4.41'
$ws.Range("E50").Value = '  -3.06%  '
$ws.Range("D51").Value = '44.07'
$ws.Range("E51").Value = '  -2.18%  '
